{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line and the\n// \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\" copyright/footer line that follow the\n// \"Requisitos\" section, collapsing the blank line left behind back down to a\n// single empty paragraph (matching the page's new, shorter footer).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two footer paragraphs by their text content rather than a fixed\n// index, so the edit still lands correctly even if earlier content shifts.\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (jupiterIndex === -1 && text.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIndex = i;\n  }\n  if (copyrightIndex === -1 && text.indexOf(\"Contact: luizeleno@usp.br\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (jupiterIndex !== -1 && copyrightIndex !== -1) {\n  // Delete the empty paragraph that sits directly above \"Ver no Jupiter ...\"\n  // too, so only one blank paragraph remains between \"Requisitos\" and the\n  // page-break paragraph that follows (same as before there were two).\n  const blankAboveIndex = jupiterIndex - 1;\n\n  // Delete from the bottom up so earlier indices stay valid.\n  items[copyrightIndex].delete();\n  items[jupiterIndex].delete();\n  if (blankAboveIndex >= 0 && items[blankAboveIndex].text.trim() === \"\") {\n    items[blankAboveIndex].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line and the\n# \"\u00a9 2020 . Contact: luizeleno@usp.br. ...\" copyright/footer line that follow the\n# \"Requisitos\" section, collapsing the blank line left behind back down to a\n# single empty paragraph (matching the page's new, shorter footer).\n\n$d = $word.ActiveDocument\n\n# Locate the two footer paragraphs by their text content rather than a fixed\n# index, so the edit still lands correctly even if earlier content shifts.\n$jupiterIndex = -1\n$copyrightIndex = -1\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($jupiterIndex -eq -1 -and $t -like \"*Ver no Jupiter*\") {\n        $jupiterIndex = $i\n    }\n    if ($copyrightIndex -eq -1 -and $t -like \"*Contact: luizeleno@usp.br*\") {\n        $copyrightIndex = $i\n    }\n    $i++\n}\n\nif ($jupiterIndex -ne -1 -and $copyrightIndex -ne -1) {\n    # Delete the empty paragraph that sits directly above \"Ver no Jupiter ...\"\n    # too, so only one blank paragraph remains between \"Requisitos\" and the\n    # page-break paragraph that follows (same as before there were two).\n    $blankAboveIndex = $jupiterIndex - 1\n\n    # Delete from the bottom up so earlier indices stay valid.\n    $d.Paragraphs.Item($copyrightIndex).Range.Delete()\n    $d.Paragraphs.Item($jupiterIndex).Range.Delete()\n    if ($blankAboveIndex -ge 1 -and $d.Paragraphs.Item($blankAboveIndex).Range.Text.Trim() -eq \"\") {\n        $d.Paragraphs.Item($blankAboveIndex).Range.Delete()\n    }\n}\n"}
